$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")
$ws.Range("C1").Value = "Web Data 13"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2"
